$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows before row 699 (existing data rows 699-749 shift to 705-755)
$ws.Range("699:704").Insert()

# Shared (unchanged) column values, copied from the template row (old row 699 / new row 705)
$mercadoId = 1
$mercado = "Agrícola del Norte S.A. de Arica"
$region = "Arica y Parinacota"
$codreg = 15
$categoriaId = 100112002
$categoria = "Pimiento"
$unidad = "`$/caja 15 kilos"
$origen = "Región de Arica y Parinacota"
$kgUnidades = 15
$clasificacion = "Hortaliza"

# New rows: row, Fecha(D), Variedad(H), Calidad(I), Volumen(J), PrecioMin(K), PrecioMax(L), PrecioProm(M), PrecioKg(P)
$newRows = @(
    @(699, 44826, "Zafiro rojo",  "Primera", 100, 12000, 13000, 12500, 833),
    @(700, 44826, "Zafiro rojo",  "Segunda", 120, 10000, 11000, 10500, 700),
    @(701, 44826, "Zafiro rojo",  "Tercera", 130, 8000,  9000,  8500,  567),
    @(702, 44826, "Zafiro verde", "Primera", 80,  8000,  9000,  8500,  567),
    @(703, 44826, "Zafiro verde", "Segunda", 120, 6000,  7000,  6500,  433),
    @(704, 44826, "Zafiro verde", "Tercera", 120, 5000,  6000,  5500,  367)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $categoriaId
    $ws.Cells.Item($r, 7).Value = $categoria
    $ws.Cells.Item($r, 8).Value = $row[2]
    $ws.Cells.Item($r, 9).Value = $row[3]
    $ws.Cells.Item($r, 10).Value = $row[4]
    $ws.Cells.Item($r, 11).Value = $row[5]
    $ws.Cells.Item($r, 12).Value = $row[6]
    $ws.Cells.Item($r, 13).Value = $row[7]
    $ws.Cells.Item($r, 14).Value = $unidad
    $ws.Cells.Item($r, 15).Value = $origen
    $ws.Cells.Item($r, 16).Value = $row[8]
    $ws.Cells.Item($r, 17).Value = $kgUnidades
    $ws.Cells.Item($r, 18).Value = $clasificacion
}
